$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6060825.5
$ws.Range("I33").Value = 213.14285
$ws.Range("J33").Value = 11363861
$ws.Range("K33").Value = 213.14285
$ws.Range("L33").Value = 11363861
$ws.Range("M33").Value = 15.85714999999999
$ws.Range("N33").Value = -11364319
$ws.Range("H92").Value = 53419570
$ws.Range("I92").Value = 2525971.8
$ws.Range("J92").Value = 333334340
$ws.Range("K92").Value = 2525971.8
$ws.Range("L92").Value = 333334340
$ws.Range("M92").Value = -2524723.8
$ws.Range("N92").Value = -333336836
$ws.Range("H98").Value = 1393.591
$ws.Range("I98").Value = 1225.5
$ws.Range("J98").Value = 2150
$ws.Range("K98").Value = 1225.5
$ws.Range("L98").Value = 2150
$ws.Range("M98").Value = 272.5
$ws.Range("N98").Value = -5146
$ws.Range("H122").Value = 1393.591
$ws.Range("I122").Value = 1225.5
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 3676.5
$ws.Range("L122").Value = 6450
$ws.Range("M122").Value = -1226.5
$ws.Range("N122").Value = -11350
$ws.Range("H129").Value = 856.05634
$ws.Range("I129").Value = 530.1429000000001
$ws.Range("J129").Value = 936.1053000000001
$ws.Range("K129").Value = 1590.4287
$ws.Range("L129").Value = 2808.3159
$ws.Range("M129").Value = 3409.5713
$ws.Range("N129").Value = -12808.3159
$ws.Range("H132").Value = 1309.9395
$ws.Range("I132").Value = 1291.2258
$ws.Range("K132").Value = 3873.6774
$ws.Range("M132").Value = -1343.6774
$ws.Range("H137").Value = 1533.1555
$ws.Range("I137").Value = 1530.8334
$ws.Range("J137").Value = 1537.8
$ws.Range("K137").Value = 4592.5002
$ws.Range("L137").Value = 4613.4
$ws.Range("M137").Value = -2042.5002
$ws.Range("N137").Value = -9713.4

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1714.4
$ws.Range("I2").Value = 1724.8334
$ws.Range("J2").Value = 1698.75
$ws.Range("K2").Value = 1724.8334
$ws.Range("L2").Value = 1698.75
$ws.Range("M2").Value = -1611.8334
$ws.Range("N2").Value = -1924.75
$ws.Range("H45").Value = 13498.5
$ws.Range("I45").Value = 17504
$ws.Range("J45").Value = 1482
$ws.Range("K45").Value = 17504
$ws.Range("L45").Value = 1482
$ws.Range("M45").Value = -17127
$ws.Range("N45").Value = -2236
$ws.Range("H61").Value = 4076.7556
$ws.Range("I61").Value = 4404.5293
$ws.Range("J61").Value = 3063.6365
$ws.Range("K61").Value = 4404.5293
$ws.Range("L61").Value = 3063.6365
$ws.Range("M61").Value = -4192.5293
$ws.Range("N61").Value = -3487.6365
$ws.Range("H74").Value = 1196
$ws.Range("I74").Value = 967.26666
$ws.Range("K74").Value = 967.26666
$ws.Range("M74").Value = -93.26666
$ws.Range("H77").Value = 1196
$ws.Range("I77").Value = 967.26666
$ws.Range("K77").Value = 4836.3333
$ws.Range("M77").Value = -468.3333000000002
$ws.Range("H110").Value = 722.8182
$ws.Range("I110").Value = 683.8
$ws.Range("K110").Value = 683.8
$ws.Range("M110").Value = 1361.2
$ws.Range("H116").Value = 1714.4
$ws.Range("I116").Value = 1724.8334
$ws.Range("J116").Value = 1698.75
$ws.Range("K116").Value = 1724.8334
$ws.Range("L116").Value = 1698.75
$ws.Range("M116").Value = 569.1666
$ws.Range("N116").Value = -6286.75
$ws.Range("H122").Value = 1711350.8
$ws.Range("I122").Value = 1833446.2
$ws.Range("K122").Value = 5500338.6
$ws.Range("M122").Value = -5497888.6
$ws.Range("H132").Value = 3422.3
$ws.Range("I132").Value = 1976.5714
$ws.Range("J132").Value = 4687.3125
$ws.Range("K132").Value = 5929.7142
$ws.Range("L132").Value = 14061.9375
$ws.Range("M132").Value = -3399.7142
$ws.Range("N132").Value = -19121.9375
$ws.Range("H136").Value = 4076.7556
$ws.Range("I136").Value = 4404.5293
$ws.Range("J136").Value = 3063.6365
$ws.Range("K136").Value = 13213.5879
$ws.Range("L136").Value = 9190.9095
$ws.Range("M136").Value = -10663.5879
$ws.Range("N136").Value = -14290.9095

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1714.4
$ws.Range("I3").Value = 1724.8334
$ws.Range("J3").Value = 1698.75
$ws.Range("K3").Value = 1724.8334
$ws.Range("L3").Value = 1698.75
$ws.Range("M3").Value = -1610.8334
$ws.Range("N3").Value = -1926.75
$ws.Range("H105").Value = 11276.652
$ws.Range("I105").Value = 23118.889
$ws.Range("J105").Value = 3663.7856
$ws.Range("K105").Value = 23118.889
$ws.Range("L105").Value = 3663.7856
$ws.Range("M105").Value = -21371.889
$ws.Range("N105").Value = -7157.7856
$ws.Range("H134").Value = 7914.35
$ws.Range("I134").Value = 12831.1
$ws.Range("K134").Value = 38493.3
$ws.Range("M134").Value = -35958.3

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 585673.2
$ws.Range("I102").Value = 942594.25
$ws.Range("K102").Value = 942594.25
$ws.Range("M102").Value = -940972.25
$ws.Range("H113").Value = 50001364
$ws.Range("I113").Value = 71429520
$ws.Range("K113").Value = 71429520
$ws.Range("M113").Value = -71427350
$ws.Range("H122").Value = 46719736
$ws.Range("I122").Value = 128473870
$ws.Range("J122").Value = 3084.0715
$ws.Range("K122").Value = 385421610
$ws.Range("L122").Value = 9252.2145
$ws.Range("M122").Value = -385419160
$ws.Range("N122").Value = -14152.2145
$ws.Range("H123").Value = 19343.121
$ws.Range("J123").Value = 19343.121
$ws.Range("L123").Value = 19343.121
$ws.Range("N123").Value = -24243.121
$ws.Range("H132").Value = 2871.25
$ws.Range("I132").Value = 2951.3333
$ws.Range("J132").Value = 2836.9285
$ws.Range("K132").Value = 8853.999899999999
$ws.Range("L132").Value = 8510.7855
$ws.Range("M132").Value = -6323.999899999999
$ws.Range("N132").Value = -13570.7855

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49114.76
$ws.Range("I7").Value = 51395.5
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 51395.5
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -51283.5
$ws.Range("N7").Value = -3724
$ws.Range("H22").Value = 644.5294
$ws.Range("I22").Value = 311.55554
$ws.Range("J22").Value = 1019.125
$ws.Range("K22").Value = 311.55554
$ws.Range("L22").Value = 1019.125
$ws.Range("M22").Value = -16.55554000000001
$ws.Range("N22").Value = -1609.125
$ws.Range("H27").Value = 644.5294
$ws.Range("I27").Value = 311.55554
$ws.Range("J27").Value = 1019.125
$ws.Range("K27").Value = 311.55554
$ws.Range("L27").Value = 1019.125
$ws.Range("M27").Value = -204.55554
$ws.Range("N27").Value = -1233.125
$ws.Range("H40").Value = 100004730
$ws.Range("I40").Value = 142859040
$ws.Range("K40").Value = 142859040
$ws.Range("M40").Value = -142858904
$ws.Range("H61").Value = 1452.5385
$ws.Range("I61").Value = 1424
$ws.Range("J61").Value = 1795
$ws.Range("K61").Value = 1424
$ws.Range("L61").Value = 1795
$ws.Range("M61").Value = -1222
$ws.Range("N61").Value = -2199
$ws.Range("H113").Value = 1452.5385
$ws.Range("I113").Value = 1424
$ws.Range("J113").Value = 1795
$ws.Range("K113").Value = 1424
$ws.Range("L113").Value = 1795
$ws.Range("M113").Value = 746
$ws.Range("N113").Value = -6135
$ws.Range("H122").Value = 3018621.8
$ws.Range("I122").Value = 4466680.5
$ws.Range("K122").Value = 13400041.5
$ws.Range("M122").Value = -13397591.5
$ws.Range("H126").Value = 49114.76
$ws.Range("I126").Value = 51395.5
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 154186.5
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -151716.5
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 9557593
$ws.Range("I132").Value = 12864667
$ws.Range("J132").Value = 3822.111
$ws.Range("K132").Value = 38594001
$ws.Range("L132").Value = 11466.333
$ws.Range("M132").Value = -38591471
$ws.Range("N132").Value = -16526.333
$ws.Range("H136").Value = 7910.878
$ws.Range("I136").Value = 5469.567
$ws.Range("K136").Value = 16408.701
$ws.Range("M136").Value = -13858.701

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1324.75
$ws.Range("I122").Value = 1263.8334
$ws.Range("J122").Value = 1507.5
$ws.Range("K122").Value = 3791.5002
$ws.Range("L122").Value = 4522.5
$ws.Range("M122").Value = -1341.5002
$ws.Range("N122").Value = -9422.5
$ws.Range("H123").Value = 29370
$ws.Range("J123").Value = 29370
$ws.Range("L123").Value = 29370
$ws.Range("N123").Value = -39170
$ws.Range("H136").Value = 2065.1042
$ws.Range("I136").Value = 2141.9644
$ws.Range("J136").Value = 1957.5
$ws.Range("K136").Value = 6425.8932
$ws.Range("L136").Value = 5872.5
$ws.Range("M136").Value = -3875.8932
$ws.Range("N136").Value = -10972.5
